$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header in C1 (copy B1's format - bold/border/center-top - onto it)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows in the new order: name, B (old), C (new)
$data = @(
    @("Alpha Acciones", 58726.12, 59541.42),
    @("Alpha Mega", 94707.83, 95954.91),
    @("Alpha Mercosur", 7848.23, 9008.200000000001),
    @("Alpha Recursos Naturales", 77333.28999999999, 76967.06),
    @("Alpha renta balan global", 41878.61, 42273.13),
    @("Arpenta acciones", 9904.09, 9899.360000000001),
    @("Arpenta ex Mercosur", 11053.06, 11044.29),
    @("Delta Acciones", 94597.02, 95044.89999999999),
    @("Delta Internacional", 4294.73, 4282.96),
    @("Delta Latinoamerica", 12684.59, 12689.04),
    @("Delta Recursos Naturales", 316993.77, 317202.17),
    @("Delta Select", 557793.5699999999, 557203.15),
    @("Delta gestion V", 99614.5, 97736.96000000001),
    @("Fima Acciones", 0, 0),
    @("Fima PB Acciones", 5890.31, 6494.77),
    @("HF Acciones Argentinas", 0, 0),
    @("HF Acciones Lideres", 35588.8, 36244.61),
    @("MAF", 0, 0),
    @("Toronto Trust Multimercado", 11655.46, 11703.27),
    @("Toronto trust Argy", 20730.87, 20804.67),
    @("avg", 73064.74000000001, 73204.74000000001),
    @("total", 1461294.85, 1464094.87)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
